$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date" placeholder text: 2/22/2024 -> 2/28/2024
#    Present on the Slide Master and on every Slide Layout (the "dt" /
#    datetimeFigureOut placeholder).
# ---------------------------------------------------------------------------
function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            $found = $tr.Find("2/22/2024", 0)
            if ($found -ne $null) {
                $found.Text = "2/28/2024"
            }
        }
    }
}

Update-DateShape($p.SlideMaster)

$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DateShape($layouts.Item($l))
}

# ---------------------------------------------------------------------------
# 2) Slide 18 ("Escrita/Leitura em Arquivos"): fix last word of body text
#    "artigo" -> "arquivo."
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(18)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        $tr = $shp.TextFrame.TextRange
        $found = $tr.Find("artigo", 0)
        if ($found -ne $null) {
            $found.Text = "arquivo."
        }
    }
}
